$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original column D (plain numeric "3", unused in the new layout) is
# removed entirely; this shifts the old E/F columns (the two SE header
# strings and the estimate/se value pairs) left into D/E and keeps the
# workbook's dimension/spans/shared-strings table consistent.
$ws.Columns("D").Delete()

# Row 2 ("Forecast"): update the SE estimates.
$ws.Range("D2").Value = 0.1
$ws.Range("E2").Value = 0.02

# Row 3 ("FE") already matches the target values, nothing to change.

# Row 4: re-label from "Forecast, FE" to "FE, Disg" and update its SE value.
$ws.Range("A4").Value = "FE"
$ws.Range("B4").Value = "Disg"
$ws.Range("D4").Value = 0.14

# Row 5: re-label from "Forecast, FE, Disg" to "FE, Var" and drop column C.
$ws.Range("A5").Value = "FE"
$ws.Range("B5").Value = "Var"
$ws.Range("C5").ClearContents()

# Row 6: re-label from "Forecast, FE, Disg" to "FE, Disg, Var" and update its SE value.
$ws.Range("A6").Value = "FE"
$ws.Range("B6").Value = "Disg"
$ws.Range("C6").Value = "Var"
$ws.Range("D6").Value = 0.14
